# Weekly refresh of the "Coco" (Vega Central Mapocho de Santiago) price
# series: two new daily observations were collected, so two new records are
# inserted at the top of their respective blocks and every older record
# shifts down one row (the two oldest records fall through to the
# newly-created rows 27 and 28 at the bottom of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New record #1: insert above row 9, pushing old rows 9-26 to 10-27 ---
$ws.Rows(9).Insert()

$ws.Cells.Item(9, 1).Value  = 9                                     # Mercado ID
$ws.Cells.Item(9, 2).Value  = "Vega Central Mapocho de Santiago"    # Mercado
$ws.Cells.Item(9, 3).Value  = "Metropolitana"                       # Región
$ws.Cells.Item(9, 4).Value  = 44424                                 # Fecha (2021-08-16)
$ws.Cells.Item(9, 5).Value  = 13                                    # Codreg
$ws.Cells.Item(9, 6).Value  = "Fruta"                                # Tipo
$ws.Cells.Item(9, 7).Value  = 100108                                # Producto ID
$ws.Cells.Item(9, 8).Value  = "Tropicales y subtropicales"          # Producto
$ws.Cells.Item(9, 9).Value  = 100108007                             # Categoría ID
$ws.Cells.Item(9, 10).Value = "Coco"                                # Categoría
$ws.Cells.Item(9, 11).Value = "Sin especificar"                     # Variedad
$ws.Cells.Item(9, 12).Value = "Primera"                             # Calidad
$ws.Cells.Item(9, 13).Value = 70                                    # Volumen
$ws.Cells.Item(9, 14).Value = 24000                                 # Precio mínimo
$ws.Cells.Item(9, 15).Value = 25000                                 # Precio máximo
$ws.Cells.Item(9, 16).Value = 24429                                 # Precio promedio ponderado
$ws.Cells.Item(9, 17).Value = "$/malla 20 unidades"                 # Unidad de comercialización
$ws.Cells.Item(9, 18).Value = "Perú"                                # Origen
$ws.Cells.Item(9, 19).Value = 1221                                  # Precio $/Kg
$ws.Cells.Item(9, 20).Value = 20                                    # Kg / unidad

# --- New record #2: insert above row 21, pushing the remaining old rows down ---
$ws.Rows(21).Insert()

$ws.Cells.Item(21, 1).Value  = 9
$ws.Cells.Item(21, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(21, 3).Value  = "Metropolitana"
$ws.Cells.Item(21, 4).Value  = 44417                                # Fecha (2021-08-09)
$ws.Cells.Item(21, 5).Value  = 13
$ws.Cells.Item(21, 6).Value  = "Fruta"
$ws.Cells.Item(21, 7).Value  = 100108
$ws.Cells.Item(21, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(21, 9).Value  = 100108007
$ws.Cells.Item(21, 10).Value = "Coco"
$ws.Cells.Item(21, 11).Value = "Sin especificar"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 30                                   # Volumen
$ws.Cells.Item(21, 14).Value = 24000                                # Precio mínimo
$ws.Cells.Item(21, 15).Value = 24000                                # Precio máximo
$ws.Cells.Item(21, 16).Value = 24000                                # Precio promedio ponderado
$ws.Cells.Item(21, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(21, 18).Value = "Perú"
$ws.Cells.Item(21, 19).Value = 1200                                 # Precio $/Kg
$ws.Cells.Item(21, 20).Value = 20
